$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text format so numeric-looking values (e.g. "1.004")
# are written as literal text and not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Apply updated values scraped on Sun Apr 23 08:22:14 UTC 2023
$ws.Range("D2").Value = '27.786.00'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '1.878.84'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '332.36'
$ws.Range("E5").Value = '  +2.65%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.4723'
$ws.Range("E7").Value = '  +4.41%  '
$ws.Range("D8").Value = '0.3956'
$ws.Range("D9").Value = '47.82'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '0.08065'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = '1.033'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '22.26'
$ws.Range("E12").Value = '  +4.21%  '
$ws.Range("D13").Value = '1.884.00'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = '5.980'
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").Value = '7.146'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").Value = '1.007'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").Value = '0.00001052'
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").Value = '87.27'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '0.06677'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = '17.24'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '27.803.53'
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").Value = '5.540'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").Value = '2.304'
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").Value = '2.117.23'
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").Value = '159.48'
$ws.Range("E27").Value = '  +3.70%  '
$ws.Range("D28").Value = '20.25'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '2.109'
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").Value = '5.621'
$ws.Range("E30").Value = '  +3.54%  '
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '0.9876'
$ws.Range("E32").Value = '  +5.77%  '
$ws.Range("D33").Value = '0.09559'
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = '1.450'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").Value = '3.593'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '5.376'
$ws.Range("E36").Value = '  +2.19%  '
$ws.Range("D37").Value = '0.06133'
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("E38").Value = '  +1.75%  '
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("D40").Value = '8.183'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("D41").Value = '0.6045'
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("D42").Value = '0.1908'
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").Value = '10.31'
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("B44").Value = 'WEMIXTOKEN'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '1.257'
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.5718'
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").Value = '12.25'
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '1.951'
$ws.Range("E47").Value = '  +1.94%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '0.06905'
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("D50").Value = '114.02'
$ws.Range("E50").Value = '  +4.94%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.00000000307'
$ws.Range("E51").Value = '  +10.92%  '

# Restore the original (default) style on the Price column now that the text values are set
$priceRange.Style = "Normal"
